# The mapping doc had a new "status" row (with a conditional lookup) added
# to the rainbows sheet, as part of testing the transform_casrec details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "rainbow"
$ws.Range("B10").Value = "status"
$ws.Range("C10").Value = "str"
$ws.Range("G10").Value = "colours"
$ws.Range("I10").Value = "purple"
$ws.Range("K10").Value = "status(blue)"
$ws.Range("J10").Value = "conditional_lookup"
$ws.Range("N10").Value = "YES"

$ws.Range("J11").Select()
